$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status text changed from "Ready for handoff" to "Handed back: in sync with en-US"
# everywhere it appears, so that the shared string itself is updated in place.
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: Latest Target File / Latest Handback File / Latest Handback DateTime
$wsZhCn.Range("I2").Value = "a.md"
$wsZhCn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-21 02:42:15"
$wsZhCn.Range("I3").Value = "a.md"
$wsZhCn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-21 02:42:15"

# de-de sheet: Latest Target File / Latest Handback File / Latest Handback DateTime
$wsDeDe.Range("I2").Value = "a.md"
$wsDeDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-21 02:42:21"
$wsDeDe.Range("I3").Value = "a.md"
$wsDeDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-21 02:42:21"

# Hyperlinks: re-add all hyperlinks for both sheets so that the new "Latest Target
# File" cells (I2/I3) get hyperlinks pointing at a.md alongside the existing A2/A3 links.
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8679fb351821d979c0657ee890c7c772ef34033/e2e/"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $baseUrl + "a.md", "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $baseUrl + "a.md", "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $baseUrl + "b.md", "", "", "b.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $baseUrl + "a.md", "", "", "a.md")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $baseUrl + "a.md", "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $baseUrl + "a.md", "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $baseUrl + "b.md", "", "", "b.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $baseUrl + "a.md", "", "", "a.md")

# Column width changes (Overview E/F, zh-cn/de-de col C and col J widened)
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527
$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZhCn.Columns.Item(10).ColumnWidth = 40
$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDeDe.Columns.Item(10).ColumnWidth = 40
